# Append a new response record (ID 7) as row 8 of Table1 on the active sheet.
# This mirrors the source workbook's existing rows (2-7), extending the
# ListObject by one row and filling in the new respondent's answers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# Grow the table by one row - this updates the table's ref/autoFilter range
# (A1:AN7 -> A1:AN8) and the sheet dimension.
$null = $lo.ListRows.Add()

# Seed the new row from the previous data row so every column keeps the same
# cell formatting (date/time format on B:C, text quote-prefix on K & Y:AM)
# as the rest of the table.
$ws.Range("A7:AN7").Copy($ws.Range("A8:AN8"))

# Now overwrite just the cells whose values differ for this new record.
$ws.Range("A8").Value = 7
$ws.Range("B8").Value = 45056.8403125
$ws.Range("C8").Value = 45056.841412037
$ws.Range("M8").Value = "MSc"
$ws.Range("W8").Value = 4
$ws.Range("X8").Value = "Sarcastic"

# These columns hold numeric-looking ratings that are stored as text
# (quote-prefixed) in the source data, so keep them as text here too.
$ws.Range("Y8").Value = "'4"
$ws.Range("Z8").Value = "'4"
$ws.Range("AA8").Value = "'4"
$ws.Range("AB8").Value = "'4"
$ws.Range("AC8").Value = "'3"
$ws.Range("AD8").Value = "'3"
$ws.Range("AE8").Value = "'5"
$ws.Range("AF8").Value = "'3"
$ws.Range("AG8").Value = "'4"
$ws.Range("AH8").Value = "'4"
$ws.Range("AI8").Value = "'4"
$ws.Range("AJ8").Value = "'3"
$ws.Range("AK8").Value = "'3"
$ws.Range("AL8").Value = "'4"
$ws.Range("AM8").Value = "'4"
